# refactor: move components that are used in Code Analysis
#
# Backlog sheet updates:
#   - "track variable in loop" (row 4) gets a completion date + Category "AST".
#   - "socket.io send pythoncode" (row 5) is renamed to "separate modules" and
#     gets a Descriptions note; a blank spacer row is inserted below it (the
#     blank Done/Descriptions cells that used to sit on row 5 slide down to
#     the new row 6), pushing the Database/Credential rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ("track variable in loop"): add completion date + category ---
# Copy the date-formatted style from D3 so D4 matches the other Done-date
# cells (numFmt "m/d/yyyy") instead of minting a brand-new style.
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = 45423
$ws.Range("E4").Value = "AST"

# --- Make room: insert a new blank row at 6, shifting the rows below down ---
$ws.Rows.Item(6).Insert()
# Insert copies row 5's formatting into the new row; row 5 only had a
# value in column B, so drop the stray inherited B6 and fix the row height
# (row-insert doesn't carry over the explicit row height).
$ws.Range("B6").Clear()
$ws.Rows.Item(6).RowHeight = 19.95

# --- Row 5 ("socket.io send pythoncode" -> "separate modules") ---
$ws.Range("B5").Value = "separate modules"
# The blank Done-date cell that used to live at D5 moved down to D6 as part
# of the row insert above, so D5 no longer holds anything.
$ws.Range("D5").Clear()
# Clear first so the new text picks up column E's default style instead of
# inheriting the old blank date-style formatting.
$ws.Range("E5").Clear()
$ws.Range("E5").Value = "separate python running modules from editor.vue"

# --- Cursor / selection left on E5 ---
$ws.Range("E5").Select()
